$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used as temporary holding space during the 3-way cell swap.
# Placed well outside the used range (A1:K18) so it has no lasting effect once cleared.
$scratch = $ws.Cells.Item(200, 26)

function Swap-RowRange($ws, $row1, $row2, $firstCol, $lastCol, $scratch) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $empty1 = ($cell1.Text -eq "")
        $empty2 = ($cell2.Text -eq "")

        if (-not $empty1) {
            $cell1.Copy($scratch)
        }

        if ($empty2) {
            $cell1.Clear()
        } else {
            $cell2.Copy($cell1)
        }

        if ($empty1) {
            $cell2.Clear()
        } else {
            $scratch.Copy($cell2)
            $scratch.Clear()
        }
    }
}

# Rows are 1-indexed in the sheet; data rows start at row 2 (row 1 is the header).
# Columns: B=2 (No.) through K=11 (bbref url); column A (index 1) is untouched.

# Swap row 6 (Tre Mann) and row 7 (Shai Gilgeous-Alexander)
Swap-RowRange $ws 6 7 2 11 $scratch

# Swap row 12 (Aleksej Pokusevski) and row 13 (Lindy Waters III)
Swap-RowRange $ws 12 13 2 11 $scratch
